$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column B (stored width 73 <- COM ColumnWidth bakes in a fixed
# ~5/6-character padding offset that the host re-applies on save, so we
# dial the input back by that amount to land exactly on width 73 in the
# saved XML, matching the wider text that now lives in column B).
$ws.Columns.Item(2).ColumnWidth = 72.16666666666667

# Turn A60:A64 into a single "Health" category cell, matching the other
# category groups already merged in column A (A53:A57, A58:A59, A2:A8,
# ...): copy the centered formatting from one of those existing merged
# category cells, merge the range, then set the label text.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A60:A64").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A60:A64").Merge() | Out-Null
$ws.Range("A60").Value = "Health"

# Scroll the view down and select the newly merged A60:A64 cell, matching
# the saved workbook's view/selection state.
$win = $excel.ActiveWindow
$win.ScrollRow = 44
$win.ScrollColumn = 1
$ws.Range("A60:A64").Select() | Out-Null
